$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.525.56'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.914.91'
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.95'
$ws.Range("E5").Value = '  +1.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4871'
$ws.Range("E7").Value = '  +3.86%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2891'
$ws.Range("E8").Value = '  +1.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06704'
$ws.Range("E9").Value = '  -3.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.02'
$ws.Range("E10").Value = '  +5.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '107.12'
$ws.Range("E11").Value = '  +1.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.920.98'
$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07604'
$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.256'
$ws.Range("E14").Value = '  +1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6642'
$ws.Range("E15").Value = '  +1.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '269.55'
$ws.Range("E16").Value = '  -5.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.502.14'
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007522'
$ws.Range("E19").Value = '  -1.14%  '

$ws.Range("E20").Value = '  -1.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.164.38'
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.489'
$ws.Range("E22").Value = '  +5.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.386'
$ws.Range("E24").Value = '  +3.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.385'
$ws.Range("E25").Value = '  +1.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.70'
$ws.Range("E26").Value = '  -2.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.98'
$ws.Range("E27").Value = '  -4.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.094'
$ws.Range("E28").Value = '  +3.25%  '

$ws.Range("E29").Value = '  -1.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.402'
$ws.Range("E30").Value = '  +2.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.114'
$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.031'
$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04983'
$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7258'
$ws.Range("E34").Value = '  -1.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").Value = '  -0.78%  '

$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02026'
$ws.Range("E38").Value = '  +1.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.674'
$ws.Range("E39").Value = '  +0.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '110.55'
$ws.Range("E40").Value = '  +2.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.005'
$ws.Range("E41").Value = '  -2.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4405'
$ws.Range("E42").Value = '  +5.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8648'
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.855'
$ws.Range("E44").Value = '  +0.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.61'
$ws.Range("E46").Value = '  +0.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.330'
$ws.Range("E47").Value = '  +2.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.282'
$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1243'
$ws.Range("E49").Value = '  +3.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '47.12'
$ws.Range("E50").Value = '  -9.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.459'
$ws.Range("E51").Value = '  +6.33%  '

